$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("products")

$ws.Range("A4").Value = "productName"
$ws.Range("B4").Value = "Bisleri"

$ws.Activate()
$ws.Range("B4").Select()
